$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Baseline formatting - copy row 8 formatting down to row 9 ---
$ws.Range("A8:AU8").Copy()
$ws.Range("A9:AU9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 2: Fix formatting for cells whose style differs from row 8 (sourced from row 7/8 cells that already use the target look) ---
$ws.Range("AS8").Copy()
$ws.Range("J9:M9").PasteSpecial(-4122)
$ws.Range("AS8").Copy()
$ws.Range("Q9:U9").PasteSpecial(-4122)
$ws.Range("AS8").Copy()
$ws.Range("Z9:AA9").PasteSpecial(-4122)
$ws.Range("AS8").Copy()
$ws.Range("AE9:AH9").PasteSpecial(-4122)
$ws.Range("AS8").Copy()
$ws.Range("AN9:AQ9").PasteSpecial(-4122)
$ws.Range("AS8").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("AK8").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("AE7").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("AS7").Copy()
$ws.Range("AR9").PasteSpecial(-4122)
$ws.Range("AI7").Copy()
$ws.Range("G9:H9").PasteSpecial(-4122)
$ws.Range("AI7").Copy()
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("AI7").Copy()
$ws.Range("W9:Y9").PasteSpecial(-4122)
$ws.Range("AI7").Copy()
$ws.Range("AB9:AD9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 3: Set the cell values/text for the new row ---
$ws.Range("A9").Value = "AB-02"
$ws.Range("B9").Value = "E4"
$ws.Range("C9").Value = "EEE-02"
$ws.Range("D9").Value = "AB-2-T9"
$ws.Range("E9").Value = 60
$ws.Range("F9").Value = "LEISURE"
$ws.Range("G9").Value = "EPI"
$ws.Range("H9").Value = "FE&HV"
$ws.Range("I9").Value = "PROJECT"
$ws.Range("J9").Value = "LEISURE"
$ws.Range("K9").Value = "LEISURE"
$ws.Range("L9").Value = "LEISURE"
$ws.Range("M9").Value = "LEISURE"
$ws.Range("N9").Value = "ACS"
$ws.Range("O9").Value = "ACS"
$ws.Range("P9").Value = "PROJECT"
$ws.Range("Q9").Value = "LEISURE"
$ws.Range("R9").Value = "LEISURE"
$ws.Range("S9").Value = "LEISURE"
$ws.Range("T9").Value = "LEISURE"
$ws.Range("U9").Value = "EPI"
$ws.Range("V9").Value = "FE&HV"
$ws.Range("W9").Value = "PROJECT"
$ws.Range("X9").Value = "ACS"
$ws.Range("Y9").Value = "ACS"
$ws.Range("Z9").Value = "LEISURE"
$ws.Range("AA9").Value = "LEISURE"
$ws.Range("AB9").Value = "HVDC&FACTS"
$ws.Range("AC9").Value = "HVDC&FACTS"
$ws.Range("AD9").Value = "PROJECT"
$ws.Range("AE9").Value = "LEISURE"
$ws.Range("AF9").Value = "LEISURE"
$ws.Range("AG9").Value = "LEISURE"
$ws.Range("AH9").Value = "LEISURE"
$ws.Range("AI9").Value = "FE&HV"
$ws.Range("AJ9").Value = "PROJECT"
$ws.Range("AK9").Value = "PROJECT"
$ws.Range("AL9").Value = "EPI"
$ws.Range("AM9").Value = "HVDC&FACTS"
$ws.Range("AN9").Value = "LEISURE"
$ws.Range("AO9").Value = "LEISURE"
$ws.Range("AP9").Value = "LEISURE"
$ws.Range("AQ9").Value = "LEISURE"
$ws.Range("AR9").Value = "PROJECT"
$ws.Range("AS9").Value = "LEISURE"
$ws.Range("AT9").Value = "LEISURE"
$ws.Range("AU9").Value = "LEISURE"

# --- Step 4: Mirror the authored selection/active cell ---
$ws.Range("D9").Select()
